$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 37, shifting existing rows 37-74 down to 38-75.
$ws.Rows.Item(37).Insert()

# Populate the newly inserted row 37 with the new daily entry.
$ws.Range("A37").Value = 5
$ws.Range("B37").Value = "Macroferia Regional de Talca"
$ws.Range("C37").Value = "Maule"
$ws.Range("D37").Value = 44587
$ws.Range("E37").Value = 7
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100103
$ws.Range("H37").Value = "Frutos de hueso (carozo)"
$ws.Range("I37").Value = 100103002
$ws.Range("J37").Value = "Ciruela"
$ws.Range("K37").Value = "Black Amber"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 200
$ws.Range("N37").Value = 9000
$ws.Range("O37").Value = 9000
$ws.Range("P37").Value = 9000
$ws.Range("Q37").Value = "$/bandeja 18 kilos granel"
$ws.Range("R37").Value = "Región de O'Higgins"
$ws.Range("S37").Value = 500
$ws.Range("T37").Value = 18
